$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 11 de Octubre de 2020 a las 21:36"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7977660
$ws.Range("C4").Value = 29370
$ws.Range("D4").Value = 5119331
$ws.Range("E4").Value = 2638729
$ws.Range("G4").Value = 230
$ws.Range("H4").Value = 219600

# Row 5 - India
$ws.Range("B5").Value = 7118770
$ws.Range("C5").Value = 67227
$ws.Range("D5").Value = 6145918
$ws.Range("E5").Value = 863673
$ws.Range("G5").Value = 808
$ws.Range("H5").Value = 109179

# Row 25 - Alemania
$ws.Range("B25").Value = 326291
$ws.Range("C25").Value = 2838
$ws.Range("E25").Value = 43089
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = 9702

# Row 27 - Israel
$ws.Range("B27").Value = 290493
$ws.Range("C27").Value = 618
$ws.Range("D27").Value = 228658
$ws.Range("E27").Value = 59855
$ws.Range("G27").Value = 39
$ws.Range("H27").Value = 1980
